# "Turret & projectile v2"
# Restructure the Bird stats table (rows 10-18): shift the header/table one
# column to the left (B..F -> A..E, keeping G as-is), and add 8 new bird
# rows (11-18) with Size/Speed-tier letters and per-row stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 10 (header row): B10->A10 (keep bold/underline style), C10->B10,
#    D10->C10, E10->D10, F10->E10 (drop), G10 unchanged.
# ---------------------------------------------------------------------

# A10 must end up with the same bold/underline style B10 currently carries.
# Copy that formatting over before B10's own content/format is overwritten.
$ws.Range("B10").Copy()
$ws.Range("A10").PasteSpecial(-4122)
[void]$excel.CutCopyMode

# B10 currently carries the bold/underline style; reset it to the default
# "Normal" style since the new B10 (Health) is unstyled.
$ws.Range("B10").Style = "Normal"

# ---------------------------------------------------------------------
# 2) New table body. Write column C (rows 11-16) first so the new shared
#    strings get created in the same order as the source edit
#    (S, M, F, XF), then A12 (Night), then D17/D18 (L, XL), then E18
#    (7+Diamond) - matching the sharedStrings table append order.
# ---------------------------------------------------------------------
$ws.Range("C11").Value = "S"
$ws.Range("C12").Value = "S"
$ws.Range("C13").Value = "M"
$ws.Range("C14").Value = "M"
$ws.Range("C15").Value = "F"
$ws.Range("C16").Value = "XF"

$ws.Range("A12").Value = "Night"

$ws.Range("D17").Value = "L"
$ws.Range("D18").Value = "XL"

$ws.Range("E18").Value = "7+Diamond"

# ---------------------------------------------------------------------
# 3) Fill in the rest of the header row.
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "Bird"
$ws.Range("B10").Value = "Health"
$ws.Range("C10").Value = "Speed"
$ws.Range("D10").Value = "Size"
$ws.Range("E10").Value = "Gold"
$ws.Range("G10").Value = "Amount (*1.1 each time up to 3x original amounts, then health x2 (&goldx2.5) until 10x original)"

# Old F10 ("Gold") is no longer part of the table - clear it.
$ws.Range("F10").ClearContents()

# ---------------------------------------------------------------------
# 4) Fill in the rest of the bird rows (11-18).
# ---------------------------------------------------------------------

# Row 11 - Flock
$ws.Range("A11").Value = "Flock"
$ws.Range("B11").Value = 3
$ws.Range("D11").Value = "M"
$ws.Range("E11").Value = 7
$ws.Range("G11").Value = 40

# Row 12 - Night
$ws.Range("B12").Value = 7
$ws.Range("D12").Value = "S"
$ws.Range("E12").Value = 8
$ws.Range("G12").Value = 30

# Row 13 - Acid
$ws.Range("A13").Value = "Acid"
$ws.Range("B13").Value = 7
$ws.Range("D13").Value = "M"
$ws.Range("E13").Value = 15
$ws.Range("G13").Value = 20

# Row 14 - Fire
$ws.Range("A14").Value = "Fire"
$ws.Range("B14").Value = 7
$ws.Range("D14").Value = "S"
$ws.Range("E14").Value = 15
$ws.Range("G14").Value = 20

# Row 15 - Thunder
$ws.Range("A15").Value = "Thunder"
$ws.Range("B15").Value = 7
$ws.Range("D15").Value = "M"
$ws.Range("E15").Value = 15
$ws.Range("G15").Value = 20

# Row 16 - Lunar
$ws.Range("A16").Value = "Lunar"
$ws.Range("B16").Value = 7
$ws.Range("D16").Value = "S"
$ws.Range("E16").Value = 50
$ws.Range("G16").Value = 10

# Row 17 - Gold
$ws.Range("A17").Value = "Gold"
$ws.Range("B17").Value = 25
$ws.Range("C17").Value = "S"
$ws.Range("E17").Value = 100
$ws.Range("G17").Value = 5

# Row 18 - Phoenix
$ws.Range("A18").Value = "Phoenix"
$ws.Range("B18").Value = 100
$ws.Range("C18").Value = "S"
$ws.Range("G18").Value = 1

# ---------------------------------------------------------------------
# 5) Styling: B11:B18 get the existing "horizontal left" style (same as
#    the rest of the workbook's left-aligned number style), E18 gets a
#    new "horizontal right" style.
# ---------------------------------------------------------------------
$ws.Range("B11:B18").HorizontalAlignment = -4131
$ws.Range("E18").HorizontalAlignment = -4152

# Row 18 needs the same explicit row height as the other table rows.
$ws.Rows.Item(18).RowHeight = 15.95

# ---------------------------------------------------------------------
# 6) Selection moves to F10 (matches the authored selection change).
# ---------------------------------------------------------------------
[void]$ws.Range("F10").Select()
